$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.043.39"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "'2.301.27"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'300.26"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").Value = "'98.17"
$ws.Range("E6").Value = "  -1.24%  "
$ws.Range("D7").Value = "'0.521"
$ws.Range("E7").Value = "  +2.40%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.516"
$ws.Range("E9").Value = "  +1.25%  "
$ws.Range("D10").Value = "'36.19"
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("D13").Value = "'17.74"
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("D14").Value = "'6.88"
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D15").Value = "'2.657.70"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").Value = "'2.277.73"
$ws.Range("E16").Value = "  -2.44%  "
$ws.Range("E17").Value = "  -1.40%  "
$ws.Range("D18").Value = "'42.935.15"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").Value = "'12.77"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").Value = "'0.0₃0911"
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D22").Value = "'68.94"
$ws.Range("E22").Value = "  +1.59%  "
$ws.Range("D23").Value = "'237.57"
$ws.Range("E23").Value = "  +0.95%  "
$ws.Range("D24").Value = "'2.14"
$ws.Range("E24").Value = "  -1.77%  "
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").Value = "'24.92"
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").Value = "'165.07"
$ws.Range("E28").Value = "  -2.63%  "
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("D31").Value = "'33.02"
$ws.Range("E31").Value = "  -3.93%  "
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("E33").Value = "  +0.70%  "
$ws.Range("E34").Value = "  +3.26%  "
$ws.Range("D35").Value = "'17.88"
$ws.Range("E35").Value = "  +1.50%  "
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("D37").Value = "'0.0698"
$ws.Range("E37").Value = "  +1.42%  "
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("E40").Value = "  -0.87%  "
$ws.Range("E41").Value = "  +0.97%  "
$ws.Range("D42").Value = "'2.015.60"
$ws.Range("E42").Value = "  +1.50%  "
$ws.Range("E43").Value = "  -1.59%  "
$ws.Range("D44").Value = "'2.21"
$ws.Range("E44").Value = "  -2.09%  "
$ws.Range("D45").Value = "'10.36"
$ws.Range("D46").Value = "'17.51"
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("D47").Value = "'2.82"
$ws.Range("E47").Value = "  -2.34%  "
$ws.Range("D48").Value = "'54.11"
$ws.Range("E48").Value = "  -2.34%  "
$ws.Range("D49").Value = "'2.525.36"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").Value = "'1.54"
$ws.Range("E50").Value = "  -0.90%  "
$ws.Range("D51").Value = "'73.13"
$ws.Range("E51").Value = "  +3.49%  "
